# ============================================================
# Applies the edits described by the commit diff to the document
# ============================================================
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $old"
    }
}

# ------------------------------------------------------------
# Hunk 1: Introduction paragraph - fix "computers how to solve"
# ------------------------------------------------------------
Replace-Text "Computer algorithms are the basis for how computers how to solve problems." "Computer algorithms are the basis for how computers solve problems."

# ------------------------------------------------------------
# Hunk 2: Introduction paragraph - reword "Visualization can explain..."
# ------------------------------------------------------------
$old2 = "Visualization can explain an algorithm" + [char]0x2019 + "s behavior when given a specific input and show the algorithm" + [char]0x2019 + "s efficiency. Having the ability to automate the visualization process can be extremely useful in algorithm development by showing designers ways in which algorithms could be changed to increase efficiency."
$new2 = "Automation of the visualization process can be even more useful in development, quickly showing ways in which algorithms could be improved to increase efficiency."
Replace-Text $old2 $new2

# ------------------------------------------------------------
# Hunk 3: Sorting Algorithms paragraph rewrite (bubble sort)
# ------------------------------------------------------------
$lq = [char]0x201c
$rq = [char]0x201d
$old3 = "One of the most studied problems in computer science is sorting. There are many different types of sorting algorithms that have been developed and each has its advantages and disadvantages. Bubble Sort is one of the most basic. It works by comparing the first 2 elements in a list and swapping them if the first is greater than the second. This repeats with the second and third elements and continues until the greatest element is ${lq}bubbled${rq} to the last position. The entire process is repeated to ${lq}bubble${rq} the second greatest element to be behind the greatest, and so on until the list is sorted."
$new3 = "Sorting is a heavily studied topic in computer science. Many types of sorting algorithms have been developed and each has its advantages and disadvantages. \textit{\textbf{Bubble Sort}} is one of the most basic. It works by comparing the first 2 elements in a list and swapping them if the first is greater than the second. This repeats with the second and third elements and so on until the greatest element is \textit{bubbled} to the last position. The process is repeated to \textit{bubble} the second greatest element to behind the greatest, and so on until the list is sorted."
Replace-Text $old3 $new3

# ------------------------------------------------------------
# Hunk 4: Decision Trees paragraph rewrite
# ------------------------------------------------------------
$old4 = "In algorithm analysis, a pruned decision tree is a tree-like structure that describes all possible execution paths the program could take, depending on the input, with any contradictory paths ${lq}pruned${rq} (removed). This decision tree is valid if there is a path from the root node to a leaf node that sorts any permutation of an `$n`$ length list. Once generated, a pruned-valid decision tree can be interpreted as the different execution paths a program could take and the efficiency of each path. The fewer the nodes in a path from the root to a leaf node, the less comparisons performed, and the more efficient the algorithm."
$new4 = "In algorithm analysis, a \textbf {pruned decision tree} is a tree that describes all possible execution paths a program can take, depending on the input, with any contradictory paths \textit{pruned}. For sorting algorithms, this tree is valid if there is a path from the root node to a leaf node that sorts any permutation of an inputted list. A \textbf{pruned-valid decision tree} can be interpreted as the different execution paths a program can take and the efficiency of each path. The fewer nodes in a path from the root to a leaf node, the less comparisons performed, and the more efficient the algorithm."
Replace-Text $old4 $new4

# ------------------------------------------------------------
# Hunk 5: Decision Tree Generator paragraph (split around the italic
# "Decision Tree Generator" run that must remain untouched)
# ------------------------------------------------------------
Replace-Text "To solve these problems, I designed an automatic analysis algorithm, the " "To solve these problems, I designed an automatic analysis code library, the "

$old5b = "This program takes a modified version of any sorting algorithm and generates a pruned-valid decision tree for some arbitrary input variables. To build the tree, the generator must run the sorting algorithm through various situations, controlling how it responds to comparisons of records."
$new5b = "This library can take a modified version of any sorting algorithm and generate a pruned-valid decision tree for some arbitrary input variables. To build the tree, the generator must run the sorting algorithm through various situations, controlling and observing its comparisons of records."
Replace-Text $old5b $new5b

# ------------------------------------------------------------
# Move the _GoBack bookmark from the end of the document to right
# after "comparisons of records." in the Decision Tree Generator
# paragraph (matches the target diff).
# ------------------------------------------------------------
$bmRange = $d.Content
$bmRange.Find.Execute("comparisons of records.")
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
